# Updated cryptos list on Mon Jun  3 21:50:59 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row in the sheet, matching a new data pull from the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Price (D) / Volume(1h) (E) text. A $null entry
# means that column is unchanged for that row.
$data = @{
    2  = @{ D = "69.137.82";  E = "  +1.88%  " }
    3  = @{ D = "3.774.51";   E = "  -0.32%  " }
    4  = @{ D = "1.00";       E = "  -0.21%  " }
    5  = @{ D = "625.68";     E = "  +3.82%  " }
    6  = @{ D = "165.58";     E = "  +1.51%  " }
    7  = @{ D = "3.769.53";   E = "  -0.42%  " }
    8  = @{ D = $null;        E = "  -0.11%  " }
    9  = @{ D = $null;        E = "  +1.55%  " }
    10 = @{ D = $null;        E = "  +1.25%  " }
    11 = @{ D = $null;        E = "  +2.99%  " }
    12 = @{ D = "6.75";       E = "  -1.32%  " }
    13 = @{ D = $null;        E = "  +0.32%  " }
    14 = @{ D = "35.61";      E = "  +1.49%  " }
    15 = @{ D = "4.403.55";   E = "  -0.41%  " }
    16 = @{ D = "3.755.10";   E = "  -0.73%  " }
    17 = @{ D = "69.119.40";  E = "  +1.89%  " }
    18 = @{ D = "17.67";      E = "  -2.77%  " }
    19 = @{ D = $null;        E = "  -1.09%  " }
    20 = @{ D = "7.05";       E = "  +0.60%  " }
    21 = @{ D = "467.09";     E = "  +2.10%  " }
    22 = @{ D = "9.57";       E = "  +1.34%  " }
    23 = @{ D = $null;        E = "  +2.26%  " }
    24 = @{ D = $null;        E = "  +3.05%  " }
    25 = @{ D = $null;        E = "  +0.23%  " }
    26 = @{ D = "12.03";      E = "  +1.70%  " }
    27 = @{ D = $null;        E = "  +3.56%  " }
    28 = @{ D = "10.02";      E = "  +1.42%  " }
    29 = @{ D = $null;        E = "  -0.07%  " }
    30 = @{ D = "3.920.86";   E = "  -0.40%  " }
    31 = @{ D = $null;        E = "  +2.62%  " }
    32 = @{ D = "2.23";       E = "  +2.41%  " }
    33 = @{ D = $null;        E = "  -0.72%  " }
    34 = @{ D = "28.78";      E = "  -0.50%  " }
    35 = @{ D = "0.175";      E = "  +21.11%  " }
    36 = @{ D = $null;        E = "  +0.27%  " }
    37 = @{ D = "3.722.84";   E = "  -0.34%  " }
    38 = @{ D = "8.95";       E = "  +0.37%  " }
    39 = @{ D = $null;        E = "  +2.16%  " }
    40 = @{ D = $null;        E = "  +5.30%  " }
    41 = @{ D = $null;        E = "  +0.46%  " }
    42 = @{ D = $null;        E = "  -0.86%  " }
    43 = @{ D = $null;        E = "  -0.11%  " }
    45 = @{ D = "153.82";     E = "  +1.16%  " }
    46 = @{ D = "43.20";      E = "  -1.06%  " }
    47 = @{ D = $null;        E = "  +0.68%  " }
    48 = @{ D = "46.68";      E = "  -0.90%  " }
    49 = @{ D = $null;        E = "  +4.07%  " }
    50 = @{ D = "8.40";       E = "  +1.51%  " }
    51 = @{ D = $null;        E = "  +0.39%  " }
}

foreach ($row in $data.Keys) {
    $entry = $data[$row]

    if ($null -ne $entry.D) {
        # Price column: some values (e.g. "1.00") look numeric and would
        # otherwise be auto-converted to a number by Excel. Force the cell
        # to text before writing, then restore the default "Normal" style
        # so the cell keeps no explicit style index, matching the rest of
        # the sheet's plain inline-string cells.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $entry.D
        $cell.Style = "Normal"
    }

    if ($null -ne $entry.E) {
        $ws.Range("E$row").Value = $entry.E
    }
}
